$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.292.24"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.864.46"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'237.05"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.4661"
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").Value = "'0.2813"
$ws.Range("E8").Value = "  +2.16%  "
$ws.Range("D9").Value = "'0.06379"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "'18.21"
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("D11").Value = "1.865.23"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'0.07570"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").Value = "'95.50"
$ws.Range("E13").Value = "  +12.66%  "
$ws.Range("D14").Value = "'4.971"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "'0.6388"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "'295.13"
$ws.Range("E16").Value = "  +22.18%  "
$ws.Range("D17").Value = "30.238.13"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "'1.005"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007398"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.111.30"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'1.008"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "'4.997"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'5.998"
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").Value = "'165.18"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "'9.026"
$ws.Range("E26").Value = "  -3.24%  "
$ws.Range("D27").Value = "'19.37"
$ws.Range("E27").Value = "  +6.94%  "
$ws.Range("D28").Value = "'1.926"
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("D29").Value = "'0.1081"
$ws.Range("E29").Value = "  +6.51%  "
$ws.Range("D30").Value = "'1.339"
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D31").Value = "'4.045"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'3.801"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "'0.04939"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").Value = "'0.7270"
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").Value = "'1.121"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("D36").Value = "'2.717"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "'0.01920"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("D38").Value = "'2.686"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'1.970"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "'0.8641"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").Value = "'106.05"
$ws.Range("E41").Value = "  +0.98%  "
$ws.Range("D42").Value = "'1.006"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "'5.575"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "'0.4057"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "'65.26"
$ws.Range("E45").Value = "  +3.24%  "
$ws.Range("D46").Value = "'7.094"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").Value = "'8.996"
$ws.Range("E47").Value = "  +4.84%  "
$ws.Range("D48").Value = "'0.1187"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").Value = "'0.05569"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'0.3719"
$ws.Range("E51").Value = "  +0.81%  "
